$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "K" values (was Strike#) computed/regenerated for rows 2-25, column G
$newK = @(7, 5, 8, 6, 0, 8, 0, 2, 6, 5, 3, 1, 3, 7, 5, 5, 3, 6, 3, 6, 4, 5, 4, 2)

for ($i = 0; $i -lt $newK.Length; $i++) {
    $row = $i + 2
    $ws.Range("G$row").Value = $newK[$i]
}
